$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "ro_CCM-code_FLASH"
$wb.Worksheets.Item(2).Name = "ro_CCM-code_CCM"
$wb.Worksheets.Item(3).Name = "ro_FLASH-code_FLASH"
$wb.Worksheets.Item(4).Name = "ro_FLASH-code_CCM"
